$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Device:" label and its "${device.deviceName}" value from row 4
# (keep the cell styling/formatting, only drop the text so the now-unused
# "Device:" / "${device.deviceName}" shared strings are dropped on save)
$ws.Range("A4:B4").ClearContents()

# Reflect the new selection (Name Box) on the sheet, now pointing at the
# cleared A4:B4 range instead of the old A10
$ws.Range("A4:B4").Select()
